$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs | Fgf18 | Fgfr2 | ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf18"
$ws.Cells.Item(2,3).Value = "Fgfr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.5972743333333334
$ws.Cells.Item(2,8).Value = 1.791823
$ws.Cells.Item(2,9).Value = 0.0994998030631086
$ws.Cells.Item(2,10).Value = 0.09949980306310859
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.8155003333333334
$ws.Cells.Item(2,14).Value = 2.446501
$ws.Cells.Item(2,15).Value = 0.1910612426590028
$ws.Cells.Item(2,16).Value = 0.1910612426590029
$ws.Cells.Item(2,17).Value = 0.4870774179247778
$ws.Cells.Item(2,18).Value = 4.383696761323
$ws.Cells.Item(2,19).Value = 0.01901055601756359
$ws.Cells.Item(2,20).Value = 0.01901055601756359

# Row 3: ECs | Fgf18 | Fgfr2 | FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf18"
$ws.Cells.Item(3,3).Value = "Fgfr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.5972743333333334
$ws.Cells.Item(3,8).Value = 1.791823
$ws.Cells.Item(3,9).Value = 0.0994998030631086
$ws.Cells.Item(3,10).Value = 0.09949980306310859
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.333134333333334
$ws.Cells.Item(3,14).Value = 9.999403000000001
$ws.Cells.Item(3,15).Value = 0.7809105179307759
$ws.Cells.Item(3,16).Value = 0.780910517930776
$ws.Cells.Item(3,17).Value = 1.990795586852111
$ws.Cells.Item(3,18).Value = 17.917160281669
$ws.Cells.Item(3,19).Value = 0.07770044274402234
$ws.Cells.Item(3,20).Value = 0.07770044274402234

# Row 4: ECs | Fgf18 | Fgfr2 | MuSCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf18"
$ws.Cells.Item(4,3).Value = "Fgfr2"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.5972743333333334
$ws.Cells.Item(4,8).Value = 1.791823
$ws.Cells.Item(4,9).Value = 0.0994998030631086
$ws.Cells.Item(4,10).Value = 0.09949980306310859
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.119632
$ws.Cells.Item(4,14).Value = 0.358896
$ws.Cells.Item(4,15).Value = 0.02802823941022116
$ws.Cells.Item(4,16).Value = 0.02802823941022117
$ws.Cells.Item(4,17).Value = 0.07145312304533334
$ws.Cells.Item(4,18).Value = 0.643078107408
$ws.Cells.Item(4,19).Value = 0.002788804301522665
$ws.Cells.Item(4,20).Value = 0.002788804301522665

# Row 5: FAPs | Fgf18 | Fgfr2 | ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf18"
$ws.Cells.Item(5,3).Value = "Fgfr2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 4.044312666666666
$ws.Cells.Item(5,8).Value = 12.132938
$ws.Cells.Item(5,9).Value = 0.6737411795567456
$ws.Cells.Item(5,10).Value = 0.6737411795567455
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.8155003333333334
$ws.Cells.Item(5,14).Value = 2.446501
$ws.Cells.Item(5,15).Value = 0.1910612426590028
$ws.Cells.Item(5,16).Value = 0.1910612426590029
$ws.Cells.Item(5,17).Value = 3.298138327770889
$ws.Cells.Item(5,18).Value = 29.683244949938
$ws.Cells.Item(5,19).Value = 0.1287258269966542
$ws.Cells.Item(5,20).Value = 0.1287258269966542

# Row 6: FAPs | Fgf18 | Fgfr2 | FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf18"
$ws.Cells.Item(6,3).Value = "Fgfr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 4.044312666666666
$ws.Cells.Item(6,8).Value = 12.132938
$ws.Cells.Item(6,9).Value = 0.6737411795567456
$ws.Cells.Item(6,10).Value = 0.6737411795567455
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.333134333333334
$ws.Cells.Item(6,14).Value = 9.999403000000001
$ws.Cells.Item(6,15).Value = 0.7809105179307759
$ws.Cells.Item(6,16).Value = 0.780910517930776
$ws.Cells.Item(6,17).Value = 13.48023740400156
$ws.Cells.Item(6,18).Value = 121.322136636014
$ws.Cells.Item(6,19).Value = 0.5261315734789501
$ws.Cells.Item(6,20).Value = 0.5261315734789501

# Row 7: FAPs | Fgf18 | Fgfr2 | MuSCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf18"
$ws.Cells.Item(7,3).Value = "Fgfr2"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 4.044312666666666
$ws.Cells.Item(7,8).Value = 12.132938
$ws.Cells.Item(7,9).Value = 0.6737411795567456
$ws.Cells.Item(7,10).Value = 0.6737411795567455
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.119632
$ws.Cells.Item(7,14).Value = 0.358896
$ws.Cells.Item(7,15).Value = 0.02802823941022116
$ws.Cells.Item(7,16).Value = 0.02802823941022117
$ws.Cells.Item(7,17).Value = 0.4838292129386667
$ws.Cells.Item(7,18).Value = 4.354462916448
$ws.Cells.Item(7,19).Value = 0.01888377908114127
$ws.Cells.Item(7,20).Value = 0.01888377908114127

# Row 8: Inflammatory-Mac | Fgf18 | Fgfr2 | ECs
$ws.Cells.Item(8,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,2).Value = "Fgf18"
$ws.Cells.Item(8,3).Value = "Fgfr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.748897
$ws.Cells.Item(8,8).Value = 2.246691
$ws.Cells.Item(8,9).Value = 0.1247585905771153
$ws.Cells.Item(8,10).Value = 0.1247585905771153
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.8155003333333334
$ws.Cells.Item(8,14).Value = 2.446501
$ws.Cells.Item(8,15).Value = 0.1910612426590028
$ws.Cells.Item(8,16).Value = 0.1910612426590029
$ws.Cells.Item(8,17).Value = 0.6107257531323333
$ws.Cells.Item(8,18).Value = 5.496531778191001
$ws.Cells.Item(8,19).Value = 0.02383653134804942
$ws.Cells.Item(8,20).Value = 0.02383653134804942

# Row 9: Inflammatory-Mac | Fgf18 | Fgfr2 | FAPs
$ws.Cells.Item(9,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,2).Value = "Fgf18"
$ws.Cells.Item(9,3).Value = "Fgfr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.748897
$ws.Cells.Item(9,8).Value = 2.246691
$ws.Cells.Item(9,9).Value = 0.1247585905771153
$ws.Cells.Item(9,10).Value = 0.1247585905771153
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.333134333333334
$ws.Cells.Item(9,14).Value = 9.999403000000001
$ws.Cells.Item(9,15).Value = 0.7809105179307759
$ws.Cells.Item(9,16).Value = 0.780910517930776
$ws.Cells.Item(9,17).Value = 2.496174302830334
$ws.Cells.Item(9,18).Value = 22.46556872547301
$ws.Cells.Item(9,19).Value = 0.09742529558388877
$ws.Cells.Item(9,20).Value = 0.09742529558388877

# Row 10: Inflammatory-Mac | Fgf18 | Fgfr2 | MuSCs
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Fgf18"
$ws.Cells.Item(10,3).Value = "Fgfr2"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.748897
$ws.Cells.Item(10,8).Value = 2.246691
$ws.Cells.Item(10,9).Value = 0.1247585905771153
$ws.Cells.Item(10,10).Value = 0.1247585905771153
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.119632
$ws.Cells.Item(10,14).Value = 0.358896
$ws.Cells.Item(10,15).Value = 0.02802823941022116
$ws.Cells.Item(10,16).Value = 0.02802823941022117
$ws.Cells.Item(10,17).Value = 0.08959204590400001
$ws.Cells.Item(10,18).Value = 0.806328413136
$ws.Cells.Item(10,19).Value = 0.003496763645177151
$ws.Cells.Item(10,20).Value = 0.003496763645177151

# Row 11: MuSCs | Fgf18 | Fgfr2 | ECs
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Fgf18"
$ws.Cells.Item(11,3).Value = "Fgfr2"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.4994553333333333
$ws.Cells.Item(11,8).Value = 1.498366
$ws.Cells.Item(11,9).Value = 0.08320415683717519
$ws.Cells.Item(11,10).Value = 0.08320415683717518
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.8155003333333334
$ws.Cells.Item(11,14).Value = 2.446501
$ws.Cells.Item(11,15).Value = 0.1910612426590028
$ws.Cells.Item(11,16).Value = 0.1910612426590029
$ws.Cells.Item(11,17).Value = 0.4073059908184444
$ws.Cells.Item(11,18).Value = 3.665753917366
$ws.Cells.Item(11,19).Value = 0.01589708959970526
$ws.Cells.Item(11,20).Value = 0.01589708959970526

# Row 12: MuSCs | Fgf18 | Fgfr2 | FAPs
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Fgf18"
$ws.Cells.Item(12,3).Value = "Fgfr2"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.4994553333333333
$ws.Cells.Item(12,8).Value = 1.498366
$ws.Cells.Item(12,9).Value = 0.08320415683717519
$ws.Cells.Item(12,10).Value = 0.08320415683717518
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 3.333134333333334
$ws.Cells.Item(12,14).Value = 9.999403000000001
$ws.Cells.Item(12,15).Value = 0.7809105179307759
$ws.Cells.Item(12,16).Value = 0.780910517930776
$ws.Cells.Item(12,17).Value = 1.664751719499778
$ws.Cells.Item(12,18).Value = 14.982765475498
$ws.Cells.Item(12,19).Value = 0.06497500120971199
$ws.Cells.Item(12,20).Value = 0.06497500120971199

# Row 13: MuSCs | Fgf18 | Fgfr2 | MuSCs
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Fgf18"
$ws.Cells.Item(13,3).Value = "Fgfr2"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.4994553333333333
$ws.Cells.Item(13,8).Value = 1.498366
$ws.Cells.Item(13,9).Value = 0.08320415683717519
$ws.Cells.Item(13,10).Value = 0.08320415683717518
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.119632
$ws.Cells.Item(13,14).Value = 0.358896
$ws.Cells.Item(13,15).Value = 0.02802823941022116
$ws.Cells.Item(13,16).Value = 0.02802823941022117
$ws.Cells.Item(13,17).Value = 0.05975084043733333
$ws.Cells.Item(13,18).Value = 0.5377575639359999
$ws.Cells.Item(13,19).Value = 0.002332066027757936
$ws.Cells.Item(13,20).Value = 0.002332066027757936

# Row 14: Resolving-Mac | Fgf18 | Fgfr2 | ECs
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Fgf18"
$ws.Cells.Item(14,3).Value = "Fgfr2"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.1128296666666667
$ws.Cells.Item(14,8).Value = 0.338489
$ws.Cells.Item(14,9).Value = 0.0187962699658552
$ws.Cells.Item(14,10).Value = 0.0187962699658552
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.8155003333333334
$ws.Cells.Item(14,14).Value = 2.446501
$ws.Cells.Item(14,15).Value = 0.1910612426590028
$ws.Cells.Item(14,16).Value = 0.1910612426590029
$ws.Cells.Item(14,17).Value = 0.09201263077655555
$ws.Cells.Item(14,18).Value = 0.828113676989
$ws.Cells.Item(14,19).Value = 0.003591238697030388
$ws.Cells.Item(14,20).Value = 0.003591238697030388

# Row 15: Resolving-Mac | Fgf18 | Fgfr2 | FAPs
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Fgf18"
$ws.Cells.Item(15,3).Value = "Fgfr2"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.1128296666666667
$ws.Cells.Item(15,8).Value = 0.338489
$ws.Cells.Item(15,9).Value = 0.0187962699658552
$ws.Cells.Item(15,10).Value = 0.0187962699658552
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 3.333134333333334
$ws.Cells.Item(15,14).Value = 9.999403000000001
$ws.Cells.Item(15,15).Value = 0.7809105179307759
$ws.Cells.Item(15,16).Value = 0.780910517930776
$ws.Cells.Item(15,17).Value = 0.3760764357852223
$ws.Cells.Item(15,18).Value = 3.384687922067
$ws.Cells.Item(15,19).Value = 0.01467820491420267
$ws.Cells.Item(15,20).Value = 0.01467820491420267

# Row 16: Resolving-Mac | Fgf18 | Fgfr2 | MuSCs
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Fgf18"
$ws.Cells.Item(16,3).Value = "Fgfr2"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.1128296666666667
$ws.Cells.Item(16,8).Value = 0.338489
$ws.Cells.Item(16,9).Value = 0.0187962699658552
$ws.Cells.Item(16,10).Value = 0.0187962699658552
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.119632
$ws.Cells.Item(16,14).Value = 0.358896
$ws.Cells.Item(16,15).Value = 0.02802823941022116
$ws.Cells.Item(16,16).Value = 0.02802823941022117
$ws.Cells.Item(16,17).Value = 0.01349803868266667
$ws.Cells.Item(16,18).Value = 0.121482348144
$ws.Cells.Item(16,19).Value = 0.0005268263546221391
$ws.Cells.Item(16,20).Value = 0.0005268263546221391

# Remove the now-unused 17th row (previously Resolving-Mac -> Resolving-Mac)
$ws.Rows.Item(17).Delete()

